$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new "2022-Q4" row at the top of the
#    data (row 2), pushing the existing rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 5 doesn't exist yet - borrow the style already used by row 4 (index
# style "2": bold / bordered / centered) so the new last row matches the
# look of the other index cells without inventing a new style entry.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

# Clear old values (formatting of A2:A4/B1:D1 is left untouched) then
# rewrite the five data rows with the quarter added at the front.
$summary.Range("A2:D5").ClearContents()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.01

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.08

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.04

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q4" detail sheet right after "总计". Duplicating the
#    existing "2022-Q3" sheet (rather than Worksheets.Add()) carries over
#    all of its number/border/alignment formatting for free.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $summary)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only had one data row (row 2); extend the index-column style
# down to rows 3-5 for the extra funds before filling in the values.
$q4.Range("A2").Copy()
$q4.Range("A3:A5").PasteSpecial(-4122)

# Force text storage for the numeric-looking string columns (fund code,
# size, position, weight, market value) so they don't get silently
# re-typed as numbers when assigned below.
$q4.Range("B2:G5").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "519097"
$q4.Range("C2").Value = "新华中小市值优选混合"
$q4.Range("D2").Value = "0.66"
$q4.Range("E2").Value = "70.51"
$q4.Range("F2").Value = "2.79"
$q4.Range("G2").Value = "0.0184"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "004250"
$q4.Range("C3").Value = "银河量化优选混合"
$q4.Range("D3").Value = "0.30"
$q4.Range("E3").Value = "62.88"
$q4.Range("F3").Value = "1.33"
$q4.Range("G3").Value = "0.0040"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "970073"
$q4.Range("C4").Value = "东证融汇成长优选混合A"
$q4.Range("D4").Value = "0.38"
$q4.Range("E4").Value = "89.59"
$q4.Range("F4").Value = "0.82"
$q4.Range("G4").Value = "0.0031"
$q4.Range("H4").Value = 7

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "970074"
$q4.Range("C5").Value = "东证融汇成长优选混合C"
$q4.Range("D5").Value = "0.11"
$q4.Range("E5").Value = "89.59"
$q4.Range("F5").Value = "0.82"
$q4.Range("G5").Value = "0.0009"
$q4.Range("H5").Value = 7
